# Updates cryptos list price/volume figures (scraper refresh).
# Note: several "Price" cells are plain decimal-looking text (e.g. "211.04")
# that Excel would otherwise auto-convert to a number on assignment; those
# are written with a leading apostrophe to force literal text, matching the
# original inline-string storage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.624.97"
$ws.Range("E2").Value = "  -1.96%  "
$ws.Range("D3").Value = "1.587.82"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'211.04"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").Value = "'19.56"
$ws.Range("E10").Value = "  -3.74%  "
$ws.Range("D11").Value = "'0.0833"
$ws.Range("E11").Value = "  -1.51%  "
$ws.Range("D12").Value = "1.810.65"
$ws.Range("E12").Value = "  -2.19%  "
$ws.Range("D13").Value = "1.596.88"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("D16").Value = "'64.79"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "26.607.64"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "0.0₃0727"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").Value = "'209.33"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'6.70"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "'4.25"
$ws.Range("E22").Value = "  -2.67%  "
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("D24").Value = "'8.85"
$ws.Range("E24").Value = "  -2.19%  "
$ws.Range("D25").Value = "'145.69"
$ws.Range("E25").Value = "  -1.66%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'7.22"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").Value = "'15.29"
$ws.Range("E29").Value = "  -2.00%  "
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("E32").Value = "  -3.97%  "
$ws.Range("D33").Value = "'0.681"
$ws.Range("E33").Value = "  +22.81%  "
$ws.Range("E34").Value = "  -2.94%  "
$ws.Range("D35").Value = "1.313.49"
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  -0.69%  "
$ws.Range("E37").Value = "  -5.02%  "
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("D39").Value = "'0.826"
$ws.Range("E39").Value = "  -3.92%  "
$ws.Range("E40").Value = "  +0.14%  "
$ws.Range("E41").Value = "  -1.34%  "
$ws.Range("D42").Value = "'5.38"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "'62.67"
$ws.Range("E44").Value = "  -4.52%  "
$ws.Range("D45").Value = "1.723.56"
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("D46").Value = "'89.26"
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "'1.60"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").Value = "'0.840"
$ws.Range("E48").Value = "  -9.59%  "
$ws.Range("E49").Value = "  -1.57%  "
$ws.Range("D50").Value = "'0.0978"
$ws.Range("E50").Value = "  -2.15%  "
$ws.Range("E51").Value = "  -0.88%  "
